$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column L (12th column).
#    This shifts old columns L..T (12-20) to M..U (13-21) for every row,
#    including both the header row and the data row, and shifts the
#    stored column-width metadata the same way.
$ws.Columns.Item(12).Insert()

# 2) Row 1 header: after the insert, the old "SIDEBAR_SUBMENU" header
#    (that used to live in L1) has moved to M1, and L1 is blank.
#    Move it back to L1, and put the new sub-submenu header in M1.
$ws.Range("M1").Cut($ws.Range("L1"))
$ws.Range("M1").Value = "SIDEBAR_SUBMENU_SUBMENU"

# 3) Row 2 data: L2 is blank after the insert (its old content "Setup
#    Jenis Parameter" moved to M2, which is correct there). Fill L2 with
#    the new value first, then apply the same alignment/number-format
#    style as K2 (set the value before pasting the format, so the paste
#    doesn't clobber it back to the default style).
$ws.Range("L2").Value = "Setup Kelengkapan Kepesertaan"
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

# 4) The newly inserted column L keeps a plain custom width corresponding
#    to a stored (OOXML) width of 15 -- Excel's ColumnWidth property is
#    offset from the stored width by 5/6, so subtract that out.
$ws.Columns.Item(12).ColumnWidth = 15 - 5/6

# 5) Restore the selection shown in the sheet view.
$ws.Range("M15").Select()
